$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = -20.94715446399935
$ws.Range("C2").Value = -24.1584550006699
$ws.Range("D2").Value = -17.59988028469123
$ws.Range("E2").Value = -0.2350536263201248
$ws.Range("F2").Value = 0.02115833168474134
$ws.Range("G2").Value = -11.10927032539325
$ws.Range("H2").Value = 0

# Row 3
$ws.Range("B3").Value = 45.55835335766614
$ws.Range("C3").Value = 42.83101988146516
$ws.Range("D3").Value = 48.33776479211864
$ws.Range("E3").Value = 0.3754068742273728
$ws.Range("F3").Value = 0.009650412162992001
$ws.Range("G3").Value = 38.90060526813625
$ws.Range("H3").Value = [double]"2.238293274943599E-306"

# Row 4
$ws.Range("B4").Value = 19.13263361458139
$ws.Range("C4").Value = 16.94218550882367
$ws.Range("D4").Value = 21.36411107927525
$ws.Range("E4").Value = 0.175067254643216
$ws.Range("F4").Value = 0.009468249444719185
$ws.Range("G4").Value = 18.48992843559459
$ws.Range("H4").Value = [double]"5.797234086427238E-75"

# Row 5
$ws.Range("B5").Value = 1.348224772020568
$ws.Range("C5").Value = -0.4238625664754525
$ws.Range("D5").Value = 3.151848717741923
$ws.Range("E5").Value = 0.01339217094023933
$ws.Range("F5").Value = 0.008999900532330248
$ws.Range("G5").Value = 1.488035439072996
$ws.Range("H5").Value = 0.1367761636417744

# Row 6
$ws.Range("B6").Value = 52.63639691012187
$ws.Range("C6").Value = 47.16294467357032
$ws.Range("D6").Value = 58.31342403063804
$ws.Range("E6").Value = 0.4228884162847065
$ws.Range("F6").Value = 0.01863171541991478
$ws.Range("G6").Value = 22.69723462138629
$ws.Range("H6").Value = [double]"5.488189412898144E-111"

# Row 7
$ws.Range("B7").Value = -2.026947786620448
$ws.Range("C7").Value = -6.609465386771952
$ws.Range("D7").Value = 2.780426300782857
$ws.Range("E7").Value = -0.02047772254523509
$ws.Range("F7").Value = 0.02444003357874808
$ws.Range("G7").Value = -0.8378762033715687
$ws.Range("H7").Value = 0.4021222186931896

# Row 8
$ws.Range("B8").Value = 8.038657326513277
$ws.Range("C8").Value = 3.199295576334427
$ws.Range("D8").Value = 13.10495301086587
$ws.Range("E8").Value = 0.07731891529998634
$ws.Range("F8").Value = 0.02338116024215527
$ws.Range("G8").Value = 3.306889585427139
$ws.Range("H8").Value = 0.0009470329022164
